# Applies the diff: renames classifier-target headers on "Valeurs réelles"
# sheet and replaces the predicted price values on both sheets with
# XGBClassifier-style integer class predictions.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Valeurs réelles" -------------------------------------------
$ws1 = $wb.Worksheets.Item("Valeurs réelles")

$ws1.Range("C1").Value = "PRIX EXP POMME GALA FRANCE 170/220G CAT.I PLATEAU 1RG_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME GALA FRANCE 170/220G CAT.I PLATEAU 1RG_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME GALA FRANCE 170/220G CAT.I PLATEAU 1RG_S+3_class"

$sheet1Data = @{
    2  = @{ C = 3; D = 2; E = 3 }
    3  = @{ C = 2; D = 2; E = 2 }
    4  = @{ C = 2; D = 2; E = 2 }
    5  = @{ C = 2; D = 2; E = 4 }
    6  = @{ C = 2; D = 4; E = 0 }
    7  = @{ C = 4; D = 0; E = 0 }
    8  = @{ C = 0; D = 0; E = 2 }
    9  = @{ C = 0; D = 2; E = 1 }
    10 = @{ C = 2; D = 1; E = 4 }
    11 = @{ C = 1; D = 4; E = 2 }
    12 = @{ C = 4; D = 2; E = 2 }
    13 = @{ C = 2; D = 2; E = 2 }
    14 = @{ C = 2; D = 2; E = 1 }
    15 = @{ C = 2; D = 1; E = 2 }
    16 = @{ C = 1; D = 2; E = 2 }
    17 = @{ C = 2; D = 2; E = 1 }
    18 = @{ C = 2; D = 1; E = 2 }
    19 = @{ C = 1; D = 2; E = 2 }
    20 = @{ C = 2; D = 2; E = 2 }
    21 = @{ C = 2; D = 2; E = 2 }
    22 = @{ C = 2; D = 2; E = 3 }
    23 = @{ C = 2; D = 3; E = 3 }
    24 = @{ C = 3; D = 3; E = 2 }
    25 = @{ C = 3; D = 2; E = 2 }
    26 = @{ C = 2; D = 2; E = 2 }
    27 = @{ C = 2; D = 2; E = 2 }
    28 = @{ C = 2; D = 2; E = 2 }
}

foreach ($row in $sheet1Data.Keys) {
    $vals = $sheet1Data[$row]
    $ws1.Range("C$row").Value = $vals.C
    $ws1.Range("D$row").Value = $vals.D
    $ws1.Range("E$row").Value = $vals.E
}

# --- Sheet 2: "Prédictions" -----------------------------------------------
$ws2 = $wb.Worksheets.Item("Prédictions")

$sheet2Data = @{
    2  = @{ B = 1;  C = 1;  D = 0 }
    3  = @{ B = 0;  C = 2;  D = 0 }
    4  = @{ B = -2; C = 2;  D = 2 }
    5  = @{ B = 2;  C = 2;  D = 2 }
    6  = @{ B = 2;  C = 2;  D = 2 }
    7  = @{ B = 2;  C = 0;  D = -2 }
    8  = @{ B = 0;  C = -1; D = 2 }
    9  = @{ B = 0;  C = -1; D = 2 }
    10 = @{ B = -2; C = 0;  D = 2 }
    11 = @{ B = -2; C = 0;  D = 2 }
    12 = @{ B = 0;  C = 1;  D = 0 }
    13 = @{ B = 0;  C = -1; D = 0 }
    14 = @{ B = -1; C = 0;  D = -2 }
    15 = @{ B = 0;  C = -1; D = 0 }
    16 = @{ B = 0;  C = 0;  D = 0 }
    17 = @{ B = 0;  C = 0;  D = -1 }
    18 = @{ B = 0;  C = 0;  D = 0 }
    19 = @{ B = 0;  C = -1; D = 0 }
    20 = @{ B = 0;  C = 0;  D = 0 }
    21 = @{ B = 0;  C = 0;  D = 0 }
    22 = @{ B = 0;  C = 0;  D = 0 }
    23 = @{ B = 0;  C = 0;  D = 0 }
    24 = @{ B = 0;  C = 0;  D = 0 }
    25 = @{ B = 0;  C = 0;  D = 0 }
    26 = @{ B = 0;  C = 0;  D = 0 }
    27 = @{ B = 0;  C = 0;  D = 0 }
    28 = @{ B = 0;  C = -1; D = 0 }
}

foreach ($row in $sheet2Data.Keys) {
    $vals = $sheet2Data[$row]
    $ws2.Range("B$row").Value = $vals.B
    $ws2.Range("C$row").Value = $vals.C
    $ws2.Range("D$row").Value = $vals.D
}
